$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.297.76"
$ws.Range("E2").Value = "  +3.77%  "
$ws.Range("D3").Value = "1.715.37"
$ws.Range("E3").Value = "  +3.30%  "
$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'239.93"
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").Value = "'0.4694"
$ws.Range("E7").Value = "  -2.08%  "
$ws.Range("D8").Value = "'0.2629"
$ws.Range("E8").Value = "  +1.24%  "
$ws.Range("D9").Value = "'0.06216"
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("D10").Value = "1.709.10"
$ws.Range("E10").Value = "  +2.86%  "
$ws.Range("D11").Value = "'0.07071"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").Value = "'15.18"
$ws.Range("E12").Value = "  +3.60%  "
$ws.Range("D13").Value = "'4.409"
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("D14").Value = "'0.5876"
$ws.Range("E14").Value = "  +0.48%  "
$ws.Range("D15").Value = "'76.19"
$ws.Range("E15").Value = "  +2.78%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").Value = "26.287.00"
$ws.Range("E18").Value = "  +3.62%  "
$ws.Range("D19").Value = "'0.000006803"
$ws.Range("E19").Value = "  +1.82%  "
$ws.Range("D20").Value = "'11.55"
$ws.Range("E20").Value = "  +1.39%  "
$ws.Range("D21").Value = "1.932.56"
$ws.Range("E21").Value = "  +3.64%  "
$ws.Range("D22").Value = "'4.553"
$ws.Range("E22").Value = "  +4.05%  "
$ws.Range("D23").Value = "'8.794"
$ws.Range("E23").Value = "  +2.45%  "
$ws.Range("D24").Value = "'5.334"
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("D25").Value = "'135.31"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("D26").Value = "'15.16"
$ws.Range("E26").Value = "  +0.60%  "
$ws.Range("D27").Value = "'1.407"
$ws.Range("E27").Value = "  +1.96%  "
$ws.Range("D28").Value = "'1.760"
$ws.Range("E28").Value = "  +5.36%  "
$ws.Range("D29").Value = "'106.85"
$ws.Range("E29").Value = "  +1.98%  "
$ws.Range("D30").Value = "'4.038"
$ws.Range("E30").Value = "  +1.17%  "
$ws.Range("D31").Value = "'3.678"
$ws.Range("E31").Value = "  +2.22%  "
$ws.Range("D32").Value = "'0.07730"
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("D33").Value = "'0.04409"
$ws.Range("E33").Value = "  +1.14%  "
$ws.Range("D34").Value = "'2.615"
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("D35").Value = "'0.6198"
$ws.Range("E35").Value = "  +3.11%  "
$ws.Range("D36").Value = "'0.9678"
$ws.Range("E36").Value = "  +3.00%  "
$ws.Range("D37").Value = "'0.9149"
$ws.Range("E37").Value = "  +7.38%  "
$ws.Range("D38").Value = "'112.95"
$ws.Range("E38").Value = "  +14.54%  "
$ws.Range("D39").Value = "'2.403"
$ws.Range("E39").Value = "  -8.43%  "
$ws.Range("D40").Value = "'1.001"
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("D41").Value = "'1.897"
$ws.Range("E41").Value = "  +4.86%  "
$ws.Range("D42").Value = "'0.01463"
$ws.Range("E42").Value = "  -2.52%  "
$ws.Range("E43").Value = "  +11.56%  "
$ws.Range("D44").Value = "'0.3799"
$ws.Range("E44").Value = "  +1.77%  "
$ws.Range("D45").Value = "'0.1145"
$ws.Range("E45").Value = "  +3.20%  "
$ws.Range("D46").Value = "'6.231"
$ws.Range("E46").Value = "  +0.92%  "
$ws.Range("E47").Value = "  +0.94%  "
$ws.Range("E48").Value = "  +4.05%  "
$ws.Range("D49").Value = "'7.675"
$ws.Range("E49").Value = "  +5.12%  "
$ws.Range("D50").Value = "'1.219"
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("D51").Value = "'0.3369"
$ws.Range("E51").Value = "  +1.37%  "
